# Edit: add "Strings" worksheet (i18n Key/Value table) between
# "Stallers and Stoppers_German" and "Evaluations", add a sheet-scoped
# defined name "list" pointing at it, and populate its data + formatting.
# Also: adding / activating the new sheet naturally clears tabSelected on
# the previously-selected first sheet and moves the workbook's activeTab.

$wb = $excel.ActiveWorkbook

$wsStoppers = $wb.Worksheets.Item("Stallers and Stoppers_German")
$newSheet = $wb.Worksheets.Add($null, $wsStoppers)
$newSheet.Name = "Strings"

$data = @(
  ,@("Key", "Value")
  ,@("Stoppers.Title", "Karrierehemmer und -stopper")
  ,@("Stoppers.Loading", "Stopper werden geladen")
  ,@("PageTitles.LIBRARY", "Kompetenzen Bibliothek")
  ,@("PageTitles.COMPETENCIES", "Kompetenzen Einschätzung")
  ,@("PageTitles.QUESTIONS", "Verhaltensbasierte Fragen")
  ,@("Home.Title", "Korn Ferry Leadership Architect")
  ,@("Home.Internal", "Nur zur internen Verwendung")
  ,@("QuestionsResult.Questions", "FRAGEN")
  ,@("QuestionsResult.Notes", "Notizen")
  ,@("QuestionsResult.Empty", "Keine Kompetenzen oder Stopper ausgewählt. Umleitung...")
  ,@("Buttons.Print", "Drucken")
  ,@("Buttons.Info", "Info")
  ,@("Buttons.Submit", "Übermitteln")
  ,@("Buttons.Reset", "Zurücksetzen")
  ,@("Buttons.Close", "Schließen")
  ,@("Questionaire.Reset", "Bist du dir sicher den Fragebogen zurückzusetzen?")
  ,@("Questionaire.PasswordRequired", "Password is required to continue")
  ,@("Questionaire.PasswordIncorrect", "Password is incorrect")
  ,@("Questionaire.Password", "Password")
  ,@("Questionaire.Login", "Login")
  ,@("Questionaire.Loading", "Kompetenzen werden geladen...")
  ,@("StopperItem.Problem", "Problem")
  ,@("StopperItem.NotAProblem", "Kein Problem")
  ,@("Library.Loading", "Kompetenzen werden geladen...")
  ,@("Library.SortByNumber", "Sortiere nach Kompetenz Nummer")
  ,@("Library.SortByFactors", "Sortiere nach Faktoren und Kategorien")
  ,@("Skills.SKILLED", "Gut ausgeprägt")
  ,@("Skills.LESS", "Schwach ausgeprägt")
  ,@("Skills.TALENTED", "Talentiert")
  ,@("Skills.OVERUSED", "Übermäßig eingesetzte Fähigkeit")
  ,@("RightsReserved", "© Korn Ferry 2014-2015. All rights reserved.")
  ,@("Evaluation.Reset", "Bist du dir sicher die Bewertung zurückzusetzen?")
  ,@("Evaluation.Loading", "Kompetenzen werden geladen...")
  ,@("CompetencyItem.Cards", "Korn Ferry Leadership Architect ©Global Competency Framework Sort Cards")
  ,@("EvaluationResult.Legend", "Legende")
  ,@("EvaluationResult.Empty", "Keine Kompetenzen bewertet. Umleitung...")
  ,@("Evaluations", "Würde beschreiben")
  ,@("Evaluations", "Könnte beschreiben")
  ,@("Evaluations", "Würde nicht beschreiben")
  ,@("Evaluations", "Dies trifft die ganze oder die meiste Zeit zu")
  ,@("Evaluations", "Dies trifft manchmal oder zeitweise zu und könnte eine Mischung aus zutreffend und nicht zutreffend sein")
  ,@("Evaluations", "Dies trifft selten oder nie zu")
)

$r = 1
foreach ($pair in $data) {
    $newSheet.Cells.Item($r, 1).Value = $pair[0]
    $newSheet.Cells.Item($r, 2).Value = $pair[1]
    $r = $r + 1
}

# Column widths approximating the source workbook (bestFit columns A & B).
$newSheet.Columns.Item(1).ColumnWidth = 30.42578125
$newSheet.Columns.Item(2).ColumnWidth = 71.140625

# Vertical-center + wrap formatting over the populated data area (matches
# the author's formatting pass recorded in the target styles).
$dataRange = $newSheet.Range("A1:B43")
$dataRange.VerticalAlignment = -4108
$dataRange.WrapText = $true

# A lingering blank formatted row (artifact of the CSV query-table refresh
# that produced this sheet) extending out to column DC.
$tailRow = $newSheet.Range("A45:DC45")
$tailRow.VerticalAlignment = -4108
$tailRow.WrapText = $true

$newSheet.Range("B7").Select()

# Sheet-scoped defined name "list" used by the app to populate dropdowns.
$newSheet.Names.Add("list", "=Strings!`$A`$2:`$B`$37")

Write-Output "Added Strings sheet with $($r - 1) rows"
